# "Pais" sheet tracks per-country COVID stats, sorted descending by column B
# (Casos totales). This refresh re-pulls live totals, which both updates
# several countries' figures and re-orders a few rows that crossed each
# other's sort position (Sudan/Libano, Gambia/Eslovenia/Lituania/Sudan del
# Sur, Santa Lucia/Timor Oriental).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 23 de Agosto de 2020 a las 22:04"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 5866780
$ws.Range("C4").Value = 25352
$ws.Range("D4").Value = 3155444
$ws.Range("E4").Value = 2530877
$ws.Range("G4").Value = 285
$ws.Range("H4").Value = 180459

# Row 6 - India
$ws.Range("B6").Value = 3105185
$ws.Range("C6").Value = 61749
$ws.Range("D6").Value = 2336796
$ws.Range("E6").Value = 710697
$ws.Range("G6").Value = 846
$ws.Range("H6").Value = 57692

# Row 23 - Alemania
$ws.Range("B23").Value = 234399
$ws.Range("C23").Value = 542
$ws.Range("E23").Value = 16117

# Row 62 - Uzbekistan
$ws.Range("B62").Value = 38946
$ws.Range("C62").Value = 534
$ws.Range("D62").Value = 34987
$ws.Range("E62").Value = 3686
$ws.Range("G62").Value = 8
$ws.Range("H62").Value = 273

# Row 86 - was Libano, now Sudan (new data pushes Sudan ahead of Libano)
$ws.Range("A86").Value = "Sudan"
$ws.Range("B86").Value = 12836
$ws.Range("C86").Value = 154
$ws.Range("D86").Value = 6497
$ws.Range("E86").Value = 5524
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 815

# Row 87 - was Sudan, now Libano (unchanged Libano figures, shifted down)
$ws.Range("A87").Value = "Libano"
$ws.Range("B87").Value = 12698
$ws.Range("C87").Value = 507
$ws.Range("D87").Value = 3625
$ws.Range("E87").Value = 8950
$ws.Range("G87").Value = 2
$ws.Range("H87").Value = 123

# Row 128 - Mali
$ws.Range("B128").Value = 2705
$ws.Range("C128").Value = 6
$ws.Range("D128").Value = 2018
$ws.Range("E128").Value = 562

# Row 129 - was Eslovenia, now Gambia (new data moves Gambia up)
$ws.Range("A129").Value = "Gambia"
$ws.Range("B129").Value = 2685
$ws.Range("C129").Value = 248
$ws.Range("D129").Value = 490
$ws.Range("E129").Value = 2108
$ws.Range("G129").Value = 3
$ws.Range("H129").Value = 87

# Row 130 - was Lituania, now Eslovenia (unchanged figures, shifted down)
$ws.Range("A130").Value = "Eslovenia"
$ws.Range("B130").Value = 2651
$ws.Range("C130").Value = 34
$ws.Range("D130").Value = 2079
$ws.Range("E130").Value = 441
$ws.Range("H130").Value = 131

# Row 131 - was Sudan del Sur, now Lituania (unchanged figures, shifted down)
$ws.Range("A131").Value = "Lituania"
$ws.Range("B131").Value = 2635
$ws.Range("C131").Value = 41
$ws.Range("D131").Value = 1766
$ws.Range("E131").Value = 785
$ws.Range("H131").Value = 84

# Row 132 - was Gambia, now Sudan del Sur (unchanged figures, shifted down)
$ws.Range("A132").Value = "Sudan del Sur"
$ws.Range("B132").Value = 2499
$ws.Range("C132").Value = 2
$ws.Range("D132").Value = 1290
$ws.Range("E132").Value = 1162
$ws.Range("H132").Value = 47

# Row 146 - Aruba
$ws.Range("B146").Value = 1568
$ws.Range("C146").Value = 34
$ws.Range("D146").Value = 455
$ws.Range("E146").Value = 1106

# Row 148 - Republica de Chipre
$ws.Range("B148").Value = 1421
$ws.Range("C148").Value = 4
$ws.Range("E148").Value = 523

# Row 202 - was Timor Oriental, now Santa Lucia (tied totals, order swaps)
$ws.Range("A202").Value = "Santa Lucia"

# Row 203 - was Santa Lucia, now Timor Oriental
$ws.Range("A203").Value = "Timor Oriental"
